$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A53").Value = "Rishi"
$ws.Range("B53").Value = 20
$ws.Range("C53").Value = "Subscribed"
$ws.Range("D53").Value = "Employed"
